$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. Paragraph 2 ("My first thought for project deliverables ...") -- split
#    the two quoted-sentence runs so that w:proofErr gramStart/gramEnd marks
#    bracket the ")" and the opening '"' the way Word's grammar checker
#    would, without changing the visible text.
# ---------------------------------------------------------------------------
$para2Xml = '<w:p ' + $wns + '>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:tab/></w:r>' +
    '<w:r><w:t xml:space="preserve">My first thought for project deliverables were the Work Breakdown Structure, however the WBS is a breakdown (no pun intended) of the whole project to completion.  </w:t></w:r>' +
    '<w:r><w:t>&quot;Project deliverables refer to the tangible or intangible outputs or outcomes that are produced as a result of completing a project.&quot; (Mathur, 2023</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>)</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:r><w:t>&quot;</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>When deliverables are defined upfront, budgeting the time, resources, and money needed to complete them is easier.&quot; (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>coAmplifi</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, n.d.)</w:t></w:r>' +
    '</w:p>'

$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.InsertXML($para2Xml)

# ---------------------------------------------------------------------------
# 2. Paragraph 3 ("As projects age, ...") -- mark "projects" with
#    gramStart/gramEnd and replace the closing sentence with the rewritten
#    conclusion about Project Deliverables.
# ---------------------------------------------------------------------------
$para3Xml = '<w:p ' + $wns + '>' +
    '<w:r><w:tab/><w:t xml:space="preserve">As projects </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">age, it can be easy to add to them. Having your deliverables set up in advance can help control your project and be used in conjunction with scope and resource management. We do this to maintain the quality of the project during the duration of the </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>projects</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> life.  &quot;Project quality focuses on the end product or service deliverables that reflect the purpose of the project&quot; (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Darnall</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, n.d.).  While </w:t></w:r>' +
    '<w:r><w:t>each</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> element </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">is </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">important </w:t></w:r>' +
    '<w:r><w:t>and can are often seen as interconnected, I believe that Project Deliverables have a distinct role in completing your project</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'

$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$r3.InsertXML($para3Xml)

# ---------------------------------------------------------------------------
# 3. Final paragraph (just a tab) -- replace with two blank paragraphs, a
#    "References:" heading and three reference entries. Each reference ends
#    with a unique placeholder run that is converted into a real hyperlink
#    afterwards (Hyperlinks.Add creates the relationship + Hyperlink style).
# ---------------------------------------------------------------------------
$para4Xml = '<w:p ' + $wns + '/>' +
    '<w:p ' + $wns + '/>' +
    '<w:p ' + $wns + '><w:r><w:t>References:</w:t></w:r></w:p>' +
    '<w:p ' + $wns + '>' +
        '<w:r><w:t xml:space="preserve">Mathur, S. (2023, August 14). Understanding project deliverables. Project Management Path. Retrieved from </w:t></w:r>' +
        '<w:r><w:t>PLACEHOLDER_LINK_1</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wns + '>' +
        '<w:r><w:br/></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Darnall</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve">, R. W., &amp; Preston, J. M. (n.d.). Project management: From simple to complex. The Open University of Hong Kong. Original source: The Saylor Foundation. Retrieved from </w:t></w:r>' +
        '<w:r><w:t>PLACEHOLDER_LINK_2</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wns + '>' +
        '<w:r><w:br/></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>coAmplifi</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve">. (n.d.). Why are project deliverables important to deadlines and success? Retrieved from </w:t></w:r>' +
        '<w:r><w:t>PLACEHOLDER_LINK_3</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wns + '/>'

$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.InsertXML($para4Xml)

# ---------------------------------------------------------------------------
# 4. Turn the three placeholder runs into real hyperlinks.
# ---------------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("PLACEHOLDER_LINK_1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($find1, "https://projectmanagementpath.com/understanding-project-deliverables/", "", "", "https://projectmanagementpath.com/understanding-project-deliverables/")

$find2 = $d.Content
$find2.Find.Execute("PLACEHOLDER_LINK_2", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($find2, "http://www.saylor.org/site/textbooks/Project%20Management%20-%20From%20Simple%20to%20Complex.pdf", "", "", "http://www.saylor.org/site/textbooks/Project%20Management%20-%20From%20Simple%20to%20Complex.pdf")

$find3 = $d.Content
$find3.Find.Execute("PLACEHOLDER_LINK_3", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($find3, "https://coamplifi.com/blog/why-are-project-deliverables-important-to-deadlines-and-success", "", "", "https://coamplifi.com/blog/why-are-project-deliverables-important-to-deadlines-and-success")

Write-Output "edit complete"
